$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 59, shifting existing rows 59:85 down to 60:86.
$ws.Rows.Item(59).Insert()

# Populate the newly inserted row 59 with the new record (same categorical
# values as the surrounding Cilantro / Macroferia Regional de Talca rows,
# with a new date, volumen prices and origin).
$ws.Cells.Item(59, 1).Value = 5
$ws.Cells.Item(59, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(59, 3).Value = "Maule"
$ws.Cells.Item(59, 4).Value = 45016
$ws.Cells.Item(59, 5).Value = 7
$ws.Cells.Item(59, 6).Value = 100112040
$ws.Cells.Item(59, 7).Value = "Cilantro"
$ws.Cells.Item(59, 8).Value = "Sin especificar"
$ws.Cells.Item(59, 9).Value = "Primera"
$ws.Cells.Item(59, 10).Value = 150
$ws.Cells.Item(59, 11).Value = 7000
$ws.Cells.Item(59, 12).Value = 7000
$ws.Cells.Item(59, 13).Value = 7000
$ws.Cells.Item(59, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(59, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(59, 16).Value = 194
$ws.Cells.Item(59, 17).Value = 36
$ws.Cells.Item(59, 18).Value = "Hortaliza"
